# Insert two new slides ("Models" and "Tableau") into the deck, right
# after "Data Manipulation" (slide 5) and before "Results" (previously
# slide 6) / "Conclusions" (previously slide 7). The two existing
# slides keep their content unchanged but shift down two positions.
#
# ppLayoutText (=2) picks the "Title and Content" layout, i.e. the
# same layout already used by the "Goals" and "Conclusions" slides
# (title placeholder + body/content placeholder).

$p = $ppt.ActivePresentation

$sModels = $p.Slides.Add(6, 2)
$sModels.Shapes.Item(1).TextFrame.TextRange.Text = "Models"

$sTableau = $p.Slides.Add(7, 2)
$sTableau.Shapes.Item(1).TextFrame.TextRange.Text = "Tableau"
